$d = $word.ActiveDocument

$oldText = "Próxima aula contexto na controler"

# Locate the paragraph that currently holds the text we need to change.
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*$oldText*") {
        $target = $p
    }
}

$rng = $target.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
  'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
  'w14:paraId="44F91665" w14:textId="1D62F1F4" w:rsidR="00090B6E" w:rsidRPr="00B41C1C" ' +
  'w:rsidRDefault="00B77971" w:rsidP="00090B6E">' +
  '<w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Próxima aula</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> seção 6</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Multiplas</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> camadas 1 aula</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$rng.InsertXML($xml)
